# Update odds values for the "Jogos da Semana" worksheet.
# This mirrors a data refresh where several odd columns were updated
# for rows 2, 3, 6 and 9 (matches in the list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 3.7
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 2.35
$ws.Range("J2").Value = 4.75
$ws.Range("L2").Value = 3.4
$ws.Range("N2").Value = 4.75
$ws.Range("X2").Value = 15
$ws.Range("AH2").Value = 9
$ws.Range("AI2").Value = 12
$ws.Range("AJ2").Value = 23
$ws.Range("AM2").Value = 5
$ws.Range("AN2").Value = 26
$ws.Range("AS2").Value = 12
$ws.Range("AV2").Value = 17

# Row 3
$ws.Range("M3").Value = 1.13
$ws.Range("N3").Value = 6
$ws.Range("Q3").Value = 2.88
$ws.Range("R3").Value = 1.4

# Row 6
$ws.Range("G6").Value = 2.6
$ws.Range("I6").Value = 2.9
$ws.Range("J6").Value = 3.4
$ws.Range("L6").Value = 3.75
$ws.Range("X6").Value = 11
$ws.Range("AA6").Value = 23
$ws.Range("AG6").Value = 7.5
$ws.Range("AH6").Value = 13
$ws.Range("AN6").Value = 15

# Row 9
$ws.Range("I9").Value = 1.3
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 2.08
$ws.Range("R9").Value = 1.73
